# Append two new inventory rows (101, 102) to the WIP tracking sheet.
# Source data is textual (dates, qty, currency, NSN stored as literal
# strings in this workbook, matching every other data row), so each
# target cell is forced to Text format before the value is written --
# otherwise Excel's COM layer would auto-coerce things like "4/23/2019"
# into a real date serial or "5" into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Row = 101
        Values = @("4/23/2019", "SPE4A4-19-V-5882", "5", '$6,637.10', "6680010876055", "TRANSMITTER,LIQUID", "GEMS", "48098", "ZZ", "2019 SEP 30")
    },
    @{
        Row = 102
        Values = @("4/23/2019", "SPE7L3-19-V-5652", "1", '$1,679.36 ', "3040012589487", "PLATE,RETAINING,SHA", "Timken", "2051G92-001", "CP", "2019 SEP 30")
    }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 10))
    $rowRange.NumberFormat = "@"

    for ($col = 1; $col -le $entry.Values.Count; $col++) {
        $ws.Cells.Item($r, $col).Value = $entry.Values[$col - 1]
    }
}
